$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain string updates (not ambiguous as pure numbers)
$ws.Range('D2').Value = '26.269.91'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.593.14'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.44%  '
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.817.36'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').Value = '1.582.84'
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('E15').Value = '  -2.52%  '
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').Value = '26.270.77'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '0.0₃0721'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  -2.45%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -1.18%  '
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').Value = '1.427.25'
$ws.Range('E33').Value = '  +6.76%  '
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('E43').Value = '  -9.86%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').Value = '1.730.32'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('E48').Value = '  -1.80%  '
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('E51').Value = '  +0.12%  '

# Updates whose new text looks like a plain number; force text storage
# so Excel keeps them as strings (matching original inline-string cells)
# instead of converting them to numeric values.
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '213.08'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.498'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '18.93'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0851'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '63.89'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.44'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '215.33'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '8.99'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '144.75'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '6.96'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.566'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.824'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.924'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.760'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '60.74'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '86.69'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0951'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.NumberFormat = "General"
$c.Style = "Normal"

